$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 data: Date, Quantity, Category(AKALI), Type(Cleaning), Amount, Description
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-10-30"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = "AKALI"
$ws.Range("D3").Value = "Cleaning"
$ws.Range("E3").Value = 200
$ws.Range("F3").Value = "test1"
